$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Cocy" -> "Cocytus" (typo fix on the showdown element name). Every cell in
# column B that held the old "Cocy" shared string gets the corrected text;
# once no cell references "Cocy" any more, the writer drops that shared
# string and appends the new "Cocytus" string at the end of the table.
$ws.Range("B15").Value = "Cocytus"
$ws.Range("B16").Value = "Cocytus"
$ws.Range("B17").Value = "Cocytus"
$ws.Range("B30").Value = "Cocytus"
$ws.Range("B31").Value = "Cocytus"

# Sheet/view fix: scroll back up near the top and leave the active
# selection on B15 (the first corrected cell) instead of the old full
# A1:A31 column selection sitting further down the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
